# testVar.docx: add support for the "for" loop tag.
#
# 1. The first paragraph ("Template de test ... variable : {var:x}") gets a
#    left tab stop at 3119 twips (= 155.95 points) in its pPr/tabs.
# 2. The trailing tab run at the end of the "Fin du gabarit" paragraph is
#    removed.
# 3. The now-redundant trailing empty paragraph is removed.

$d = $word.ActiveDocument

# --- 1. Add the tab stop to the first paragraph ---------------------------
# Word's TabStops.Add takes a position in points; 3119 dxa / 20 = 155.95 pt.
$d.Paragraphs(1).TabStops.Add(155.95)

# --- 2. Strip the trailing tab character from the "Fin du gabarit" line ---
$finParagraph = $d.Paragraphs(2)
$tabPos = $finParagraph.Range.End - 2
$d.Range($tabPos, $tabPos + 1).Delete()

# --- 3. Remove the trailing empty paragraph --------------------------------
$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)
$d.Range($lastParagraph.Range.Start - 1, $lastParagraph.Range.End).Delete()
